$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Hunk 0: ALC row 86
$ws.Cells.Item(86, 8).Value = 224738.33
$ws.Cells.Item(86, 9).Value = 3217.75
$ws.Cells.Item(86, 10).Value = 401954.8
$ws.Cells.Item(86, 11).Value = 3217.75
$ws.Cells.Item(86, 12).Value = 401954.8
$ws.Cells.Item(86, 13).Value = -2094.75
$ws.Cells.Item(86, 14).Value = -404200.8

# Hunk 1: ALC row 89
$ws.Cells.Item(89, 8).Value = 224738.33
$ws.Cells.Item(89, 9).Value = 3217.75
$ws.Cells.Item(89, 10).Value = 401954.8
$ws.Cells.Item(89, 11).Value = 16088.75
$ws.Cells.Item(89, 12).Value = 2009774
$ws.Cells.Item(89, 13).Value = -10472.75
$ws.Cells.Item(89, 14).Value = -2021006

# Hunk 2: ALC row 98
$ws.Cells.Item(98, 8).Value = 879
$ws.Cells.Item(98, 9).Value = 854.8
$ws.Cells.Item(98, 10).Value = 1000
$ws.Cells.Item(98, 11).Value = 854.8
$ws.Cells.Item(98, 12).Value = 1000
$ws.Cells.Item(98, 13).Value = 643.2
$ws.Cells.Item(98, 14).Value = -3996

# Hunk 3: ALC row 99
$ws.Cells.Item(99, 8).Value = 523.3333
$ws.Cells.Item(99, 9).Value = 228.2
$ws.Cells.Item(99, 10).Value = 1999
$ws.Cells.Item(99, 11).Value = 684.5999999999999
$ws.Cells.Item(99, 12).Value = 5997
$ws.Cells.Item(99, 13).Value = 813.4000000000001
$ws.Cells.Item(99, 14).Value = -8993

# Hunk 4: ALC row 101
$ws.Cells.Item(101, 8).Value = 25004724
$ws.Cells.Item(101, 9).Value = 50008450
$ws.Cells.Item(101, 10).Value = 997.5
$ws.Cells.Item(101, 11).Value = 150025350
$ws.Cells.Item(101, 12).Value = 2992.5
$ws.Cells.Item(101, 13).Value = -150023728
$ws.Cells.Item(101, 14).Value = -6236.5

# Hunk 5: ALC row 113
$ws.Cells.Item(113, 8).Value = 4224
$ws.Cells.Item(113, 9).Value = 4298.6
$ws.Cells.Item(113, 10).Value = 4099.6665
$ws.Cells.Item(113, 11).Value = 4298.6
$ws.Cells.Item(113, 12).Value = 4099.6665
$ws.Cells.Item(113, 13).Value = -1044.6
$ws.Cells.Item(113, 14).Value = -10607.6665

# Hunk 6: ALC row 122
$ws.Cells.Item(122, 8).Value = 879
$ws.Cells.Item(122, 9).Value = 854.8
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 2564.4
$ws.Cells.Item(122, 12).Value = 3000
$ws.Cells.Item(122, 13).Value = -114.3999999999996
$ws.Cells.Item(122, 14).Value = -7900

# Hunk 7: ALC row 137
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 8).Value = 1999
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 1999
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 5997
$ws.Cells.Item(137, 14).Value = -11097

$ws = $wb.Worksheets.Item("ARM")
# Hunk 8: ARM row 61
$ws.Cells.Item(61, 8).Value = 2161
$ws.Cells.Item(61, 9).Value = 1806.25
$ws.Cells.Item(61, 10).Value = 4999
$ws.Cells.Item(61, 11).Value = 1806.25
$ws.Cells.Item(61, 12).Value = 4999
$ws.Cells.Item(61, 13).Value = -1594.25
$ws.Cells.Item(61, 14).Value = -5423

# Hunk 9: ARM row 132
$ws.Cells.Item(132, 8).Value = 1553
$ws.Cells.Item(132, 9).Value = 1553
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 4659
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -2129

# Hunk 10: ARM row 136
$ws.Cells.Item(136, 8).Value = 2161
$ws.Cells.Item(136, 9).Value = 1806.25
$ws.Cells.Item(136, 10).Value = 4999
$ws.Cells.Item(136, 11).Value = 5418.75
$ws.Cells.Item(136, 12).Value = 14997
$ws.Cells.Item(136, 13).Value = -2868.75
$ws.Cells.Item(136, 14).Value = -20097

$ws = $wb.Worksheets.Item("BSM")
# Hunk 11: BSM row 22
$ws.Cells.Item(22, 8).Value = 200
$ws.Cells.Item(22, 9).Value = 200
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 200
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -27

# Hunk 12: BSM row 25
$ws.Cells.Item(25, 8).Value = 134
$ws.Cells.Item(25, 9).Value = 134
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 134
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 101

# Hunk 13: BSM row 107
$ws.Cells.Item(107, 8).Value = 812.9
$ws.Cells.Item(107, 9).Value = 828.75
$ws.Cells.Item(107, 10).Value = 749.5
$ws.Cells.Item(107, 11).Value = 828.75
$ws.Cells.Item(107, 12).Value = 749.5
$ws.Cells.Item(107, 13).Value = 1091.25
$ws.Cells.Item(107, 14).Value = -4589.5

$ws = $wb.Worksheets.Item("CRP")
# Hunk 14: CRP row 5
$ws.Cells.Item(5, 8).Value = 884
$ws.Cells.Item(5, 9).Value = 551.75
$ws.Cells.Item(5, 10).Value = 1216.25
$ws.Cells.Item(5, 11).Value = 551.75
$ws.Cells.Item(5, 12).Value = 1216.25
$ws.Cells.Item(5, 13).Value = -439.75
$ws.Cells.Item(5, 14).Value = -1440.25

# Hunk 15: CRP row 15
$ws.Cells.Item(15, 13).ClearContents()
$ws.Cells.Item(15, 8).Value = 750
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 750
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 750
$ws.Cells.Item(15, 14).Value = -1090

# Hunk 16: CRP row 22
$ws.Cells.Item(22, 8).Value = 396.5
$ws.Cells.Item(22, 9).Value = 396.5
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 396.5
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -46.5

# Hunk 17: CRP row 31
$ws.Cells.Item(31, 8).Value = 3141.3333
$ws.Cells.Item(31, 9).Value = 3083.25
$ws.Cells.Item(31, 10).Value = 3257.5
$ws.Cells.Item(31, 11).Value = 3083.25
$ws.Cells.Item(31, 12).Value = 3257.5
$ws.Cells.Item(31, 13).Value = -2788.25
$ws.Cells.Item(31, 14).Value = -3847.5

# Hunk 18: CRP row 34
$ws.Cells.Item(34, 8).Value = 3141.3333
$ws.Cells.Item(34, 9).Value = 3083.25
$ws.Cells.Item(34, 10).Value = 3257.5
$ws.Cells.Item(34, 11).Value = 3083.25
$ws.Cells.Item(34, 12).Value = 3257.5
$ws.Cells.Item(34, 13).Value = -2881.25
$ws.Cells.Item(34, 14).Value = -3661.5

$ws = $wb.Worksheets.Item("CUL")
# Hunk 19: CUL row 12
$ws.Cells.Item(12, 8).Value = 221.38461
$ws.Cells.Item(12, 9).Value = 352.33334
$ws.Cells.Item(12, 10).Value = 109.14286
$ws.Cells.Item(12, 11).Value = 1057.00002
$ws.Cells.Item(12, 12).Value = 327.42858
$ws.Cells.Item(12, 13).Value = -884.0000199999999
$ws.Cells.Item(12, 14).Value = -673.42858

# Hunk 20: CUL row 16
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0

# Hunk 21: CUL row 23
$ws.Cells.Item(23, 8).Value = 1327.6
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 1327.6
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 3982.8
$ws.Cells.Item(23, 14).Value = -4452.799999999999

# Hunk 22: CUL row 76
$ws.Cells.Item(76, 8).Value = 13546.637
$ws.Cells.Item(76, 9).Value = 10253.25
$ws.Cells.Item(76, 10).Value = 15428.571
$ws.Cells.Item(76, 11).Value = 30759.75
$ws.Cells.Item(76, 12).Value = 46285.713
$ws.Cells.Item(76, 13).Value = -30376.75
$ws.Cells.Item(76, 14).Value = -47051.713

# Hunk 23: CUL row 79
$ws.Cells.Item(79, 8).Value = 13546.637
$ws.Cells.Item(79, 9).Value = 10253.25
$ws.Cells.Item(79, 10).Value = 15428.571
$ws.Cells.Item(79, 11).Value = 30759.75
$ws.Cells.Item(79, 12).Value = 46285.713
$ws.Cells.Item(79, 13).Value = -29433.75
$ws.Cells.Item(79, 14).Value = -48937.713

# Hunk 24: CUL row 106
$ws.Cells.Item(106, 8).Value = 18724.875
$ws.Cells.Item(106, 9).Value = 16933.334
$ws.Cells.Item(106, 10).Value = 19799.8
$ws.Cells.Item(106, 11).Value = 50800.00199999999
$ws.Cells.Item(106, 12).Value = 59399.39999999999
$ws.Cells.Item(106, 13).Value = -49854.00199999999
$ws.Cells.Item(106, 14).Value = -61291.39999999999

# Hunk 25: CUL row 112
$ws.Cells.Item(112, 8).Value = 10759.625
$ws.Cells.Item(112, 9).Value = 1519.25
$ws.Cells.Item(112, 10).Value = 20000
$ws.Cells.Item(112, 11).Value = 4557.75
$ws.Cells.Item(112, 12).Value = 60000
$ws.Cells.Item(112, 13).Value = -3449.75
$ws.Cells.Item(112, 14).Value = -62216

# Hunk 26: CUL row 137
$ws.Cells.Item(137, 8).Value = 3337.7778
$ws.Cells.Item(137, 9).Value = 1697.5
$ws.Cells.Item(137, 10).Value = 3806.4285
$ws.Cells.Item(137, 11).Value = 5092.5
$ws.Cells.Item(137, 12).Value = 11419.2855
$ws.Cells.Item(137, 13).Value = 7.5
$ws.Cells.Item(137, 14).Value = -21619.2855

# Hunk 27: CUL row 140
$ws.Cells.Item(140, 8).Value = 1565.2858
$ws.Cells.Item(140, 9).Value = 865
$ws.Cells.Item(140, 10).Value = 2499
$ws.Cells.Item(140, 11).Value = 2595
$ws.Cells.Item(140, 12).Value = 7497
$ws.Cells.Item(140, 13).Value = 2585
$ws.Cells.Item(140, 14).Value = -17857

$ws = $wb.Worksheets.Item("GSM")
# Hunk 28: GSM row 22
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0

# Hunk 29: GSM row 122
$ws.Cells.Item(122, 8).Value = 8309
$ws.Cells.Item(122, 9).Value = 7394.8
$ws.Cells.Item(122, 10).Value = 9832.666999999999
$ws.Cells.Item(122, 11).Value = 22184.4
$ws.Cells.Item(122, 12).Value = 29498.001
$ws.Cells.Item(122, 13).Value = -19734.4
$ws.Cells.Item(122, 14).Value = -34398.001

$ws = $wb.Worksheets.Item("LTW")
# Hunk 30: LTW row 7
$ws.Cells.Item(7, 8).Value = 6590.136
$ws.Cells.Item(7, 9).Value = 5443.3335
$ws.Cells.Item(7, 10).Value = 7384.077
$ws.Cells.Item(7, 11).Value = 5443.3335
$ws.Cells.Item(7, 12).Value = 7384.077
$ws.Cells.Item(7, 13).Value = -5331.3335
$ws.Cells.Item(7, 14).Value = -7608.077

# Hunk 31: LTW row 21
$ws.Cells.Item(21, 8).Value = 3458.3333
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 3458.3333
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 3458.3333
$ws.Cells.Item(21, 14).Value = -3806.3333

# Hunk 32: LTW row 22
$ws.Cells.Item(22, 8).Value = 11
$ws.Cells.Item(22, 9).Value = 11
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 11
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 284

# Hunk 33: LTW row 27
$ws.Cells.Item(27, 8).Value = 11
$ws.Cells.Item(27, 9).Value = 11
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 11
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = 96

# Hunk 34: LTW row 40
$ws.Cells.Item(40, 8).Value = 4239.4
$ws.Cells.Item(40, 9).Value = 3232.3333
$ws.Cells.Item(40, 10).Value = 5750
$ws.Cells.Item(40, 11).Value = 3232.3333
$ws.Cells.Item(40, 12).Value = 5750
$ws.Cells.Item(40, 13).Value = -3096.3333
$ws.Cells.Item(40, 14).Value = -6022

# Hunk 35: LTW row 122
$ws.Cells.Item(122, 8).Value = 6292.2915
$ws.Cells.Item(122, 9).Value = 5030.5386
$ws.Cells.Item(122, 10).Value = 7783.4546
$ws.Cells.Item(122, 11).Value = 15091.6158
$ws.Cells.Item(122, 12).Value = 23350.3638
$ws.Cells.Item(122, 13).Value = -12641.6158
$ws.Cells.Item(122, 14).Value = -28250.3638

# Hunk 36: LTW row 126
$ws.Cells.Item(126, 8).Value = 6590.136
$ws.Cells.Item(126, 9).Value = 5443.3335
$ws.Cells.Item(126, 10).Value = 7384.077
$ws.Cells.Item(126, 11).Value = 16330.0005
$ws.Cells.Item(126, 12).Value = 22152.231
$ws.Cells.Item(126, 13).Value = -13860.0005
$ws.Cells.Item(126, 14).Value = -27092.231

# Hunk 37: LTW row 136
$ws.Cells.Item(136, 8).Value = 3339.2
$ws.Cells.Item(136, 9).Value = 4599.5
$ws.Cells.Item(136, 10).Value = 2499
$ws.Cells.Item(136, 11).Value = 13798.5
$ws.Cells.Item(136, 12).Value = 7497
$ws.Cells.Item(136, 13).Value = -11248.5
$ws.Cells.Item(136, 14).Value = -12597

$ws = $wb.Worksheets.Item("WVR")
# Hunk 38: WVR row 29
$ws.Cells.Item(29, 14).ClearContents()
$ws.Cells.Item(29, 8).Value = 449
$ws.Cells.Item(29, 9).Value = 449
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 449
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = -159

# Hunk 39: WVR row 41
$ws.Cells.Item(41, 8).Value = 17590.857
$ws.Cells.Item(41, 9).Value = 16439.666
$ws.Cells.Item(41, 10).Value = 18454.25
$ws.Cells.Item(41, 11).Value = 16439.666
$ws.Cells.Item(41, 12).Value = 18454.25
$ws.Cells.Item(41, 13).Value = -16049.666
$ws.Cells.Item(41, 14).Value = -19234.25

# Hunk 40: WVR row 132
$ws.Cells.Item(132, 8).Value = 1448
$ws.Cells.Item(132, 9).Value = 1448
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 4344
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -1814
